# Finished aggregating g_H and building counterfactuals.
# The id2g_H sheet's DK_Central block gains a new "GT" plant id (row 18,
# id_DK_Central_GT) which shifts every following DK_Central row down by
# one, and the whole DK_Decentral block (previously rows 22-42) is
# removed entirely - the last two DK_Central rows (IndustryH, SH) take
# over rows 22-23, and the sheet's used range shrinks to A1:B23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contents for rows 7 through 23 (column A = id, column B = g_H group).
$rows = @(
    @{Row=7;  A="id_DK_Central_IndustryH_Biomass"; B="DK_Central"},
    @{Row=8;  A="id_DK_Central_BP_Coal";            B="DK_Central"},
    @{Row=9;  A="id_DK_Central_BH_Natgas";          B="DK_Central"},
    @{Row=10; A="id_DK_Central_BP_Natgas";           B="DK_Central"},
    @{Row=11; A="id_DK_Central_IndustryH_Natgas";    B="DK_Central"},
    @{Row=12; A="id_DK_Central_BH_Oil";              B="DK_Central"},
    @{Row=13; A="id_DK_Central_BP_Oil";              B="DK_Central"},
    @{Row=14; A="id_DK_Central_IndustryH_Oil";       B="DK_Central"},
    @{Row=15; A="id_DK_Central_BH_Waste";            B="DK_Central"},
    @{Row=16; A="id_DK_Central_BP_Waste";            B="DK_Central"},
    @{Row=17; A="id_DK_Central_EP";                  B="DK_Central"},
    @{Row=18; A="id_DK_Central_GT";                  B="DK_Central"},
    @{Row=19; A="id_DK_Central_HPstandard";          B="DK_Central"},
    @{Row=20; A="id_DK_Central_HPsurplusheat";       B="DK_Central"},
    @{Row=21; A="id_DK_Central_IH";                  B="DK_Central"},
    @{Row=22; A="id_DK_Central_IndustryH";           B="DK_Central"},
    @{Row=23; A="id_DK_Central_SH";                  B="DK_Central"}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}

# Drop the old DK_Decentral block (previously rows 22-42); the used
# range now ends at row 23.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -ge 24) {
    $deleteRange = $ws.Range("A24:B$lastRow")
    $deleteRange.ClearContents()
}
